$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly report adds two price rows (Angeleno - Primera / Segunda,
# Región de O'Higgins) ahead of the existing data block. Insert two blank
# rows at row 29 so the previous rows 29-33 shift down to 31-35 (matching
# the dimension change from A1:T33 to A1:T35), then populate the two new
# rows with the new week's figures.
$ws.Range("A29:A30").EntireRow.Insert()

# Row 29: Ciruela / Angeleno / Primera
$ws.Range("A29").Value = 2
$ws.Range("B29").Value = "Comercializadora del Agro de Limarí"
$ws.Range("C29").Value = "Coquimbo"
$ws.Range("D29").Value = 44636
$ws.Range("E29").Value = 4
$ws.Range("F29").Value = "Fruta"
$ws.Range("G29").Value = 100103
$ws.Range("H29").Value = "Frutos de hueso (carozo)"
$ws.Range("I29").Value = 100103002
$ws.Range("J29").Value = "Ciruela"
$ws.Range("K29").Value = "Angeleno"
$ws.Range("L29").Value = "Primera"
$ws.Range("M29").Value = 16
$ws.Range("N29").Value = 235000
$ws.Range("O29").Value = 240000
$ws.Range("P29").Value = 237500
$ws.Range("Q29").Value = "$/bins (450 kilos)"
$ws.Range("R29").Value = "Región de O'Higgins"
$ws.Range("S29").Value = 528
$ws.Range("T29").Value = 450

# Row 30: Ciruela / Angeleno / Segunda
$ws.Range("A30").Value = 2
$ws.Range("B30").Value = "Comercializadora del Agro de Limarí"
$ws.Range("C30").Value = "Coquimbo"
$ws.Range("D30").Value = 44636
$ws.Range("E30").Value = 4
$ws.Range("F30").Value = "Fruta"
$ws.Range("G30").Value = 100103
$ws.Range("H30").Value = "Frutos de hueso (carozo)"
$ws.Range("I30").Value = 100103002
$ws.Range("J30").Value = "Ciruela"
$ws.Range("K30").Value = "Angeleno"
$ws.Range("L30").Value = "Segunda"
$ws.Range("M30").Value = 20
$ws.Range("N30").Value = 185000
$ws.Range("O30").Value = 190000
$ws.Range("P30").Value = 187500
$ws.Range("Q30").Value = "$/bins (450 kilos)"
$ws.Range("R30").Value = "Región de O'Higgins"
$ws.Range("S30").Value = 417
$ws.Range("T30").Value = 450
